$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect (same password hash as the
# original file) so the locked cells below can be edited, then re-protect.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer note (A33):
# 2021-06-09 -> 2021-06-10
$ws.Range("A33").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) / Percent Change (E) columns for rows 2-30 with the
# new model-holdings snapshot values.
$ws.Range("D2").Value = 0.004847356702150149
$ws.Range("E2").Value = 0.0100871160018341
$ws.Range("D3").Value = 0.01345760200373622
$ws.Range("E3").Value = 0.04049117898500043
$ws.Range("D4").Value = 0.3085120483164205
$ws.Range("E4").Value = 0.01129180959658482
$ws.Range("D5").Value = 0.3260168669762271
$ws.Range("E5").Value = 0.02087682672233826
$ws.Range("D6").Value = 0.01795031998672752
$ws.Range("E6").Value = -0.008023283253362656
$ws.Range("D7").Value = 0.001517078705713794
$ws.Range("E7").Value = 0.005860048259220818
$ws.Range("D8").Value = 0.003319688253661207
$ws.Range("E8").Value = -0.0148865784499056
$ws.Range("D9").Value = 0.00338997276867256
$ws.Range("E9").Value = 0.0186967789707515
$ws.Range("D10").Value = 0.002901328045260299
$ws.Range("E10").Value = 0.007750540735400158
$ws.Range("D11").Value = 0.003236329982137327
$ws.Range("E11").Value = 0.006043370067543519
$ws.Range("D12").Value = 0.01727043235305
$ws.Range("E12").Value = 0.006691900075700197
$ws.Range("D13").Value = 0.03294832427779588
$ws.Range("E13").Value = 0.003156113750244272
$ws.Range("D14").Value = 0.002980816484856472
$ws.Range("E14").Value = 0.006666666666666821
$ws.Range("D15").Value = 0.01558716005454894
$ws.Range("E15").Value = 0.008998127906274522
$ws.Range("D16").Value = 0.0119293321624477
$ws.Range("E16").Value = -0.0155885602062108
$ws.Range("D17").Value = 0.03705068603883611
$ws.Range("E17").Value = 0.00266974926375485
$ws.Range("D18").Value = 0.05901365566941491
$ws.Range("E18").Value = 0.0143933120391182
$ws.Range("D19").Value = 0.007584033857891666
$ws.Range("E19").Value = -0.006523058252427161
$ws.Range("D20").Value = 0.0217860032624252
$ws.Range("E20").Value = 0.003845433727478254
$ws.Range("D21").Value = 0.004419975601755596
$ws.Range("E21").Value = -0.02626597255087548
$ws.Range("D22").Value = 0.005513987546723972
$ws.Range("E22").Value = 0.02226858877086491
$ws.Range("D23").Value = 0.001353212241585763
$ws.Range("E23").Value = 0.02185380557648853
$ws.Range("D24").Value = 0.009868385185812725
$ws.Range("E24").Value = 0.007418947993174552
$ws.Range("D25").Value = 0.006173166350085653
$ws.Range("E25").Value = 0.01736625863018348
$ws.Range("D26").Value = 0.03356085591792086
$ws.Range("E26").Value = 0.0009723261032161812
$ws.Range("D27").Value = 0.003132602797966926
$ws.Range("E27").Value = 0.005083260297984449
$ws.Range("D28").Value = 0.02672704650374413
$ws.Range("E28").Value = 0.007059532521200174
$ws.Range("D29").Value = 0.01795173195243088
$ws.Range("E29").Value = 0.003010679391047599
$ws.Range("D30").Value = 0.9999999999999999
$ws.Range("E30").Value = 0.01248886279367079

$ws.Protect("D382")
